$d = $word.ActiveDocument

# --- Paragraph 1 --------------------------------------------------------
# The visible text doesn't change ("This is demo for gitlatch changes done
# text"), but the run boundaries are restructured:
#   "This is demo for " + "gitlatch"   -> merged into "This is demo for gitlatch"
#   " " + "te"                         -> merged into " te"
# Re-typing the paragraph's text (via Find/Replace) collapses it down to a
# single run; we then nudge the formatting (no-op Bold toggle) on the
# sub-ranges that must stay separate so the engine keeps them split into
# distinct <w:r> runs at the boundaries the target document expects.

$para1 = $d.Paragraphs(1).Range
$fullText = $para1.Text.Substring(0, $para1.Text.Length - 1)

$searchRange = $d.Range($para1.Start, $para1.End)
[void]$searchRange.Find.Execute($fullText, $false, $false, $false, $false, $false, `
                                 $true, 1, $false, $fullText, 2)

$seg1 = $d.Range($para1.Start + 0, $para1.Start + 25)   # "This is demo for gitlatch"
$seg1.Font.Bold = 1
$seg1.Font.Bold = 0

$seg3 = $d.Range($para1.Start + 38, $para1.Start + 41)  # " te"
$seg3.Font.Bold = 1
$seg3.Font.Bold = 0

$seg4 = $d.Range($para1.Start + 41, $para1.Start + 43)  # "xt"
$seg4.Font.Bold = 1
$seg4.Font.Bold = 0

# --- Paragraph 2 ---------------------------------------------------------
# "Done" gains a new trailing run: " test" (kept as its own run, not merged
# into the "Done" run).
$para2 = $d.Paragraphs(2).Range
$insertPoint = $d.Range($para2.End - 1, $para2.End - 1)
$insertPoint.InsertAfter(" test")

$newRun = $d.Range($para2.End - 6, $para2.End - 1)      # " test"
$newRun.Font.Bold = 1
$newRun.Font.Bold = 0
